# Auto-generated Excel COM-interop script to apply the commit diff
# Updates specific cell values across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 13261.305
$ws.Range("I62").Value = 5328.4546
$ws.Range("K62").Value = 5328.4546
$ws.Range("M62").Value = -4704.4546
$ws.Range("H65").Value = 13261.305
$ws.Range("I65").Value = 5328.4546
$ws.Range("K65").Value = 26642.273
$ws.Range("M65").Value = -23522.273
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H96").Value = 10678.667
$ws.Range("I96").Value = 10678.667
$ws.Range("K96").Value = 32036.001
$ws.Range("M96").Value = -30663.001
$ws.Range("H98").Value = 3570.4
$ws.Range("I98").Value = 2087.25
$ws.Range("K98").Value = 2087.25
$ws.Range("M98").Value = -589.25
$ws.Range("H106").Value = 8076.909
$ws.Range("I106").Value = 1399.4286
$ws.Range("K106").Value = 1399.4286
$ws.Range("M106").Value = -768.4286
$ws.Range("H122").Value = 3570.4
$ws.Range("I122").Value = 2087.25
$ws.Range("K122").Value = 6261.75
$ws.Range("M122").Value = -3811.75
$ws.Range("H138").Value = 3296.1226
$ws.Range("I138").Value = 3200
$ws.Range("J138").Value = 3298.125
$ws.Range("K138").Value = 9600
$ws.Range("L138").Value = 9894.375
$ws.Range("M138").Value = -4460
$ws.Range("N138").Value = -20174.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13858.33
$ws.Range("I32").Value = 9624.406000000001
$ws.Range("J32").Value = 56668
$ws.Range("K32").Value = 9624.406000000001
$ws.Range("L32").Value = 56668
$ws.Range("M32").Value = -9337.406000000001
$ws.Range("N32").Value = -57242
$ws.Range("H45").Value = 22728714
$ws.Range("I45").Value = 23810890
$ws.Range("K45").Value = 23810890
$ws.Range("M45").Value = -23810513
$ws.Range("H61").Value = 13165195
$ws.Range("I61").Value = 5268.3076
$ws.Range("K61").Value = 5268.3076
$ws.Range("M61").Value = -5056.3076
$ws.Range("H64").Value = 46499.668
$ws.Range("I64").Value = 44999
$ws.Range("K64").Value = 44999
$ws.Range("M64").Value = -44751
$ws.Range("H67").Value = 46499.668
$ws.Range("I67").Value = 44999
$ws.Range("K67").Value = 44999
$ws.Range("M67").Value = -44141
$ws.Range("H97").Value = 861.88464
$ws.Range("I97").Value = 840
$ws.Range("J97").Value = 982.25
$ws.Range("K97").Value = 840
$ws.Range("L97").Value = 982.25
$ws.Range("M97").Value = -344
$ws.Range("N97").Value = -1974.25
$ws.Range("H132").Value = 4239.3687
$ws.Range("I132").Value = 2714.6382
$ws.Range("K132").Value = 8143.9146
$ws.Range("M132").Value = -5613.9146
$ws.Range("H136").Value = 13165195
$ws.Range("I136").Value = 5268.3076
$ws.Range("K136").Value = 15804.9228
$ws.Range("M136").Value = -13254.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1299.7333
$ws.Range("I64").Value = 1321.8
$ws.Range("J64").Value = 1288.7
$ws.Range("K64").Value = 1321.8
$ws.Range("L64").Value = 1288.7
$ws.Range("M64").Value = -1096.8
$ws.Range("N64").Value = -1738.7
$ws.Range("H67").Value = 1299.7333
$ws.Range("I67").Value = 1321.8
$ws.Range("J67").Value = 1288.7
$ws.Range("K67").Value = 1321.8
$ws.Range("L67").Value = 1288.7
$ws.Range("M67").Value = -541.8
$ws.Range("N67").Value = -2848.7
$ws.Range("H86").Value = 3671.125
$ws.Range("I86").Value = 3941.2
$ws.Range("K86").Value = 3941.2
$ws.Range("M86").Value = -2818.2
$ws.Range("H89").Value = 3671.125
$ws.Range("I89").Value = 3941.2
$ws.Range("K89").Value = 19706
$ws.Range("M89").Value = -14090
$ws.Range("H94").Value = 581
$ws.Range("I94").Value = 473.375
$ws.Range("J94").Value = 868
$ws.Range("K94").Value = 473.375
$ws.Range("L94").Value = 868
$ws.Range("M94").Value = -22.375
$ws.Range("N94").Value = -1770
$ws.Range("H105").Value = 1855.069
$ws.Range("I105").Value = 1661.3636
$ws.Range("K105").Value = 1661.3636
$ws.Range("M105").Value = 85.63640000000009
$ws.Range("H118").Value = 48326.668
$ws.Range("J118").Value = 48326.668
$ws.Range("L118").Value = 48326.668
$ws.Range("N118").Value = -51640.668
$ws.Range("I134").Value = 1898.0857
$ws.Range("K134").Value = 5694.257100000001
$ws.Range("M134").Value = -3159.257100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15448.772
$ws.Range("I31").Value = 15741.929
$ws.Range("K31").Value = 15741.929
$ws.Range("M31").Value = -15446.929
$ws.Range("H34").Value = 15448.772
$ws.Range("I34").Value = 15741.929
$ws.Range("K34").Value = 15741.929
$ws.Range("M34").Value = -15539.929
$ws.Range("H58").Value = 4159.1816
$ws.Range("I58").Value = 2726.4
$ws.Range("K58").Value = 2726.4
$ws.Range("M58").Value = -2523.4
$ws.Range("H112").Value = 73146.125
$ws.Range("J112").Value = 73146.125
$ws.Range("L112").Value = 73146.125
$ws.Range("N112").Value = -76100.125
$ws.Range("H136").Value = 4159.1816
$ws.Range("I136").Value = 2726.4
$ws.Range("K136").Value = 8179.200000000001
$ws.Range("M136").Value = -5629.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2291.3333
$ws.Range("J58").Value = 2437
$ws.Range("L58").Value = 7311
$ws.Range("N58").Value = -7567
$ws.Range("H75").Value = 670.6667
$ws.Range("I75").Value = 670.6667
$ws.Range("K75").Value = 2012.0001
$ws.Range("M75").Value = -1014.0001
$ws.Range("H78").Value = 670.6667
$ws.Range("I78").Value = 670.6667
$ws.Range("K78").Value = 6036.0003
$ws.Range("M78").Value = -1044.0003
$ws.Range("H119").Value = 12249.875
$ws.Range("I119").Value = 3999.5
$ws.Range("K119").Value = 11998.5
$ws.Range("M119").Value = -7160.5
$ws.Range("H131").Value = 33633.332
$ws.Range("J131").Value = 37450
$ws.Range("L131").Value = 112350
$ws.Range("N131").Value = -122430

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 35000
$ws.Range("I62").Value = 35000
$ws.Range("K62").Value = 35000
$ws.Range("M62").Value = -34314
$ws.Range("H65").Value = 35000
$ws.Range("I65").Value = 35000
$ws.Range("K65").Value = 105000
$ws.Range("M65").Value = -101568
$ws.Range("H107").Value = 879.5625
$ws.Range("I107").Value = 1006
$ws.Range("J107").Value = 601.4
$ws.Range("K107").Value = 1006
$ws.Range("L107").Value = 601.4
$ws.Range("M107").Value = 914
$ws.Range("N107").Value = -4441.4
$ws.Range("H122").Value = 3553.577
$ws.Range("I122").Value = 3540.5881
$ws.Range("K122").Value = 10621.7643
$ws.Range("M122").Value = -8171.764299999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1201.9412
$ws.Range("I16").Value = 1239.5625
$ws.Range("K16").Value = 1239.5625
$ws.Range("M16").Value = -1069.5625
$ws.Range("H22").Value = 865.5
$ws.Range("I22").Value = 878.4286
$ws.Range("J22").Value = 775
$ws.Range("K22").Value = 878.4286
$ws.Range("L22").Value = 775
$ws.Range("M22").Value = -583.4286
$ws.Range("N22").Value = -1365
$ws.Range("H27").Value = 865.5
$ws.Range("I27").Value = 878.4286
$ws.Range("J27").Value = 775
$ws.Range("K27").Value = 878.4286
$ws.Range("L27").Value = 775
$ws.Range("M27").Value = -771.4286
$ws.Range("N27").Value = -989
$ws.Range("H43").Value = 58854.285
$ws.Range("I43").Value = 58000
$ws.Range("J43").Value = 59993.332
$ws.Range("K43").Value = 58000
$ws.Range("L43").Value = 59993.332
$ws.Range("M43").Value = -57807
$ws.Range("N43").Value = -60379.332
$ws.Range("H82").Value = 1548.2963
$ws.Range("J82").Value = 1697.6471
$ws.Range("L82").Value = 1697.6471
$ws.Range("N82").Value = -2419.6471
$ws.Range("H85").Value = 1548.2963
$ws.Range("J85").Value = 1697.6471
$ws.Range("L85").Value = 1697.6471
$ws.Range("N85").Value = -4193.6471
$ws.Range("H136").Value = 85214.664
$ws.Range("J136").Value = 120171.37
$ws.Range("L136").Value = 360514.11
$ws.Range("N136").Value = -365614.11

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1866.375
$ws.Range("J96").Value = 1711
$ws.Range("L96").Value = 1711
$ws.Range("N96").Value = -4457
$ws.Range("H100").Value = 953.21875
$ws.Range("I100").Value = 860.48
$ws.Range("J100").Value = 1284.4286
$ws.Range("K100").Value = 1720.96
$ws.Range("L100").Value = 2568.8572
$ws.Range("M100").Value = -1179.96
$ws.Range("N100").Value = -3650.8572
$ws.Range("H113").Value = 432.30768
$ws.Range("I113").Value = 456
$ws.Range("J113").Value = 148
$ws.Range("K113").Value = 1368
$ws.Range("L113").Value = 444
$ws.Range("M113").Value = 802
$ws.Range("N113").Value = -4784
$ws.Range("H132").Value = 2384184.5
$ws.Range("I132").Value = 3480.0334
$ws.Range("K132").Value = 10440.1002
$ws.Range("M132").Value = -7910.100199999999
$ws.Range("H136").Value = 1072304.5
$ws.Range("I136").Value = 2875.2188
$ws.Range("J136").Value = 2783391.2
$ws.Range("K136").Value = 8625.6564
$ws.Range("L136").Value = 8350173.600000001
$ws.Range("M136").Value = -6075.6564
$ws.Range("N136").Value = -8355273.600000001
